# Set the diagonal "self-to-self" cable length cells to 0.
# These are the cells where the row country equals the column country
# (e.g. B2, C3, D4, ... Z26), which previously held non-zero lengths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("B2", "C3", "D4", "E5", "H8", "I9", "K11", "L12", "M13", "N14", "P16", "Q17", "R18", "T20", "Y25", "Z26")

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 0
}
